# Apply Diebold-Mariano correction values to the P_valores and Estadisticos_DM sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: P_valores ---
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.3590997471984654
$wsP.Range("D2").Value = 0.9670935300015311
$wsP.Range("E2").Value = 0.768339686013521
$wsP.Range("F2").Value = 0.5084733096611793

$wsP.Range("B3").Value = 0.3590997471984654
$wsP.Range("D3").Value = 0.4686305257590322
$wsP.Range("E3").Value = 0.5987985341715687
$wsP.Range("F3").Value = 0.886037845656416

$wsP.Range("B4").Value = 0.9670935300015311
$wsP.Range("C4").Value = 0.4686305257590322
$wsP.Range("E4").Value = 0.5761156918454575
$wsP.Range("F4").Value = 0.3912223736152693

$wsP.Range("B5").Value = 0.768339686013521
$wsP.Range("C5").Value = 0.5987985341715687
$wsP.Range("D5").Value = 0.5761156918454575
$wsP.Range("F5").Value = 0.6321422126746556

$wsP.Range("B6").Value = 0.5084733096611793
$wsP.Range("C6").Value = 0.886037845656416
$wsP.Range("D6").Value = 0.3912223736152693
$wsP.Range("E6").Value = 0.6321422126746556

# --- Sheet: Estadisticos_DM ---
$wsE = $wb.Worksheets.Item("Estadisticos_DM")

$wsE.Range("C2").Value = 0.948218489413392
$wsE.Range("D2").Value = -0.04199781583525591
$wsE.Range("E2").Value = 0.300328930190103
$wsE.Range("F2").Value = 0.6785688267966051

$wsE.Range("B3").Value = -0.948218489413392
$wsE.Range("D3").Value = -0.7449389116060402
$wsE.Range("E3").Value = -0.5383388273454124
$wsE.Range("F3").Value = -0.1459547105358091

$wsE.Range("B4").Value = 0.04199781583525591
$wsE.Range("C4").Value = 0.7449389116060402
$wsE.Range("E4").Value = 0.5724191171045615
$wsE.Range("F4").Value = 0.8847718767600312

$wsE.Range("B5").Value = -0.300328930190103
$wsE.Range("C5").Value = 0.5383388273454124
$wsE.Range("D5").Value = -0.5724191171045615
$wsE.Range("F5").Value = 0.4893929207338083

$wsE.Range("B6").Value = -0.6785688267966051
$wsE.Range("C6").Value = 0.1459547105358091
$wsE.Range("D6").Value = -0.8847718767600312
$wsE.Range("E6").Value = -0.4893929207338083
